$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.511.23"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.640.99"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.0000"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.0000"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "303.80"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3797"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "51.84"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.88%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08183"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.237"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.57"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.468"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.368"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001240"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.90%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.634.66"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "95.25"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.26%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.56"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.570"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.52"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.498.66"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.513"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.084"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -5.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.21"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "152.49"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.258"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.34"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.817.19"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.098"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +13.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.612"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.02%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -7.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.57"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.67%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2506"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.08766"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.019"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.07068"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7063"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.75%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.32"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.72"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6557"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9992"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.288"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.968"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.74%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "128.92"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.92%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.97%  "
